$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6.087478991596672
$ws.Range("G2").Value = 6.632064627284575
$ws.Range("H2").Value = 6.532849534453788
$ws.Range("F3").Value = 7.07791596638653
$ws.Range("G3").Value = 6.934768102664178
$ws.Range("H3").Value = 6.832509762093848
$ws.Range("F4").Value = 6.629642296918772
$ws.Range("G4").Value = 7.049883377082129
$ws.Range("H4").Value = 6.94990825857364
$ws.Range("F5").Value = 7.500521848739495
$ws.Range("G5").Value = 7.164795115113377
$ws.Range("H5").Value = 7.063406280703199
$ws.Range("F6").Value = 6.677237535014005
$ws.Range("G6").Value = 7.097467445262037
$ws.Range("H6").Value = 6.99759762443464
$ws.Range("F7").Value = 7.07214306086071
$ws.Range("G7").Value = 6.736394667201086
$ws.Range("H7").Value = 6.635210343247174
$ws.Range("F8").Value = 6.927415788133472
$ws.Range("G8").Value = 6.591638873269897
$ws.Range("H8").Value = 6.490723769041995
$ws.Range("F9").Value = 6.687184443007967
$ws.Range("G9").Value = 6.796180113582604
$ws.Range("H9").Value = 6.695634697635617
$ws.Range("F10").Value = 6.982490713208363
$ws.Range("G10").Value = 6.646765413225019
$ws.Range("H10").Value = 6.545363040921147
$ws.Range("F11").Value = 6.621971988795511
$ws.Range("G11").Value = 6.73093246143
$ws.Range("H11").Value = 6.630719031321313
$ws.Range("F12").Value = 6.676951680672263
$ws.Range("G12").Value = 6.533762090407221
$ws.Range("H12").Value = 6.431898067813862
$ws.Range("F13").Value = 6.622156162464992
$ws.Range("G13").Value = 6.731168164256895
$ws.Range("H13").Value = 6.630468669190537
$ws.Range("F14").Value = 6.585105508870231
$ws.Range("G14").Value = 6.694114806625117
$ws.Range("H14").Value = 6.593440825069703
$ws.Range("F15").Value = 7.00377532065456
$ws.Range("G15").Value = 6.633364304659977
$ws.Range("H15").Value = 6.531819100672684
$ws.Range("F16").Value = 6.390982013858166
$ws.Range("G16").Value = 6.499972782710001
$ws.Range("H16").Value = 6.399473607363614
$ws.Range("F17").Value = 6.750402801120443
$ws.Range("G17").Value = 6.4146698621414
$ws.Range("H17").Value = 6.313339622504748
$ws.Range("F18").Value = 6.322084542908061
$ws.Range("G18").Value = 6.431068045664651
$ws.Range("H18").Value = 6.330637410741721
$ws.Range("F19").Value = 6.180588213150886
$ws.Range("G19").Value = 6.289582612757311
$ws.Range("H19").Value = 6.189049186779274
$ws.Range("F20").Value = 6.553782718304688
$ws.Range("G20").Value = 6.2180432363297
$ws.Range("H20").Value = 6.116774775479422
$ws.Range("F21").Value = 6.025976750700269
$ws.Range("G21").Value = 6.134936842972269
$ws.Range("H21").Value = 6.034726999725964
$ws.Range("F22").Value = 5.881292135315664
$ws.Range("G22").Value = 5.990245273527879
$ws.Range("H22").Value = 5.890101005207271
$ws.Range("F23").Value = 5.498371438147911
$ws.Range("G23").Value = 5.918590897439366
$ws.Range("H23").Value = 5.818819584832688
$ws.Range("F24").Value = 5.520908496732022
$ws.Range("G24").Value = 5.629886308814923
$ws.Range("H24").Value = 5.529509349679705
$ws.Range("F25").Value = 5.770747628706672
$ws.Range("G25").Value = 5.434974314537371
$ws.Range("H25").Value = 5.334025226916956
$ws.Range("F26").Value = 4.7951571525162
$ws.Range("G26").Value = 5.339785756949938
$ws.Range("H26").Value = 5.240165846126454
$ws.Range("F27").Value = 5.097489653926078
$ws.Range("G27").Value = 5.206463388227134
$ws.Range("H27").Value = 5.106124889630204
$ws.Range("F28").Value = 4.505415819553637
$ws.Range("G28").Value = 5.050013822030307
$ws.Range("H28").Value = 4.950682239173527
$ws.Range("F29").Value = 5.162116628469567
$ws.Range("G29").Value = 4.791717453407658
$ws.Range("H29").Value = 4.69006040949986
$ws.Range("F30").Value = 5.089661777438251
$ws.Range("G30").Value = 4.753920103582146
$ws.Range("H30").Value = 4.65267233744249
$ws.Range("F31").Value = 4.613322128851524
$ws.Range("G31").Value = 4.722303526162458
$ws.Range("H31").Value = 4.621892750695715
$ws.Range("F32").Value = 4.083701213818873
$ws.Range("G32").Value = 4.628337326109993
$ws.Range("H32").Value = 4.52864666194311
